$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    $cell.Formula = "=""2012-05-18"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
